# Apply the changes described by the diff to DH_technology_cost.xlsx
$wb = $excel.ActiveWorkbook

# --- Sheet 2 "prices and emmision factors": fix "Nan" -> "NaN" typo, move selection to C2
$ws2 = $wb.Worksheets.Item("prices and emmision factors")
$ws2.Range("B4").Value = "NaN"
$ws2.Range("C2").Select()

# --- Sheet 3 "financal and other parameteres": insert a new header row (row 2) with
#     human readable labels, shift the old data row down to row 3, and fix a data value.
$ws3 = $wb.Worksheets.Item("financal and other parameteres")
$ws3.Rows.Item(2).Insert()

# Give the new header row the same formatting as the row above it (row 1)
$ws3.Range("A1:J1").Copy()
$ws3.Range("A2:J2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new human readable header labels (typed in the same order the
# author originally entered them, so shared-string ids line up)
$ws3.Range("C2").Value = "heat_storage"
$ws3.Range("D2").Value = "Interes Rate [0-1]"
$ws3.Range("B2").Value = "Demand "
$ws3.Range("F2").Value = "Temperature"
$ws3.Range("I2").Value = "Demand Scaling Factor [0-1]"
$ws3.Range("J2").Value = "Total Demand[ MWh]"
$ws3.Range("H2").Value = "Total Renewable Factor [0-1]"
$ws3.Range("G2").Value = "Threshold Temperature [Celsius]"
$ws3.Range("E2").Value = "Radiation"
$ws3.Range("A2").Value = "CO2 Price"

# Fix data value that changed when the row shifted down (H column, was 0.5 now 0)
$ws3.Range("H3").Value = 0

# Re-fit the columns whose widest entry changed because of the new header text
# (columns B, C and E keep their original best-fit width unchanged)
$ws3.Columns.Item(1).ColumnWidth = 8.592447916666666
$ws3.Columns.Item(4).ColumnWidth = 15.592447916666666
$ws3.Columns.Item(6).ColumnWidth = 10.592447916666666
$ws3.Columns.Item(7).ColumnWidth = 21.307291666666668
$ws3.Columns.Item(8).ColumnWidth = 21.166666666666668
$ws3.Columns.Item(9).ColumnWidth = 20.022135416666668
$ws3.Columns.Item(10).ColumnWidth = 19.451822916666668

$ws3.Range("F17").Select()

$wb.Save()
